$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fermentativos = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$bebidas        = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"

# Swap the two requirement rows so that "Tecnologia de Bebidas" now comes
# before "Tecnologia de Processos Fermentativos" (matching new shared-string order).
$ws.Range("B24").Value = $bebidas
$ws.Range("C24").Value = $bebidas

$ws.Range("B25").Value = $fermentativos
$ws.Range("C25").Value = $fermentativos
